# Update the 5x3 lattice-multiplication table: each cell gets a new
# "A x B" problem, a new 2-digit multiplier line, and new leading lattice
# digits, while keeping the "----" separator line unchanged.
$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# New content for every cell, in row-major (top-left to bottom-right) order.
# Each entry is: header ("A x B"), multiplier line, first lattice digit cell,
# second lattice digit cell. The "  ----" separator line is unchanged by the
# edit and is spliced back in below. Pieces are joined with a vertical-tab
# (Word's manual line break) to reproduce the original <w:br/>-separated runs.
$newCells = @(
    @("56 x 91", "  9    1", "5|    |", "6|    |"),
    @("26 x 48", "  4    8", "2|    |", "6|    |"),
    @("86 x 83", "  8    3", "8|    |", "6|    |"),
    @("62 x 74", "  7    4", "6|    |", "2|    |"),
    @("44 x 17", "  1    7", "4|    |", "4|    |"),
    @("30 x 98", "  9    8", "3|    |", "0|    |"),
    @("54 x 71", "  7    1", "5|    |", "4|    |"),
    @("73 x 12", "  1    2", "7|    |", "3|    |"),
    @("64 x 77", "  7    7", "6|    |", "4|    |"),
    @("83 x 39", "  3    9", "8|    |", "3|    |"),
    @("90 x 98", "  9    8", "9|    |", "0|    |"),
    @("55 x 30", "  3    0", "5|    |", "5|    |"),
    @("88 x 14", "  1    4", "8|    |", "8|    |"),
    @("95 x 74", "  7    4", "9|    |", "5|    |"),
    @("84 x 35", "  3    5", "8|    |", "4|    |"),
)

$nl = [char]11   # vertical tab == Word's manual line break (<w:br/>)
$separator = "  ----"
$rows = $table.Rows.Count
$cols = $table.Columns.Count

if ($rows * $cols -ne $newCells.Count) {
    throw "Expected $($newCells.Count) cells but table is $rows x $cols."
}

$i = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $p = $newCells[$i]
        $i++
        $cell = $table.Cell($r, $c)
        $cell.Range.Text = ($p[0] + $nl + $p[1] + $nl + $separator + $nl + $p[2] + $nl + $p[3])
    }
}

Write-Host "Updated $i cells in the lattice multiplication table."
